# adding log level for debugging
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partner Match")

# Row 36: Terranova Corporation - remove its "Algramo" solver match, count back to 0
$ws.Range("B36").Value = "['None']"
$ws.Range("C36").Value = 0

# Row 38: The Pershing Square Foundation - remove its "Algramo" solver match, count back to 0
$ws.Range("B38").Value = "['None']"
$ws.Range("C38").Value = 0

# Row 40: Yum! Brands - now matched with Mycotech, count to 1
$ws.Range("B40").Value = "['None'],Mycotech"
$ws.Range("C40").Value = 1
